$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case connector words (de/del/y/el/los/las/la) in municipality/state names ---
$textChanges = @{
    "B6" = "Pabellón De Arteaga"
    "B7" = "Rincón De Romos"
    "B8" = "San Francisco De Los Romo"
    "B28" = "Amatenango De La Frontera"
    "B31" = "Bejucal De Ocampo"
    "B36" = "Chiapa De Corzo"
    "B39" = "Comitán De Domínguez"
    "B55" = "Ocozocoautla De Espinosa"
    "B60" = "San Cristóbal De Las Casas"
    "B101" = "Guadalupe Y Calvo"
    "B104" = "Hidalgo Del Parral"
    "B125" = "San Francisco De Borja"
    "B126" = "San Francisco De Conchos"
    "B127" = "San Francisco Del Oro"
    "B135" = "Valle De Zaragoza"
    "B156" = "Villa De Álvarez"
    "A158" = "Ciudad De México"
    "B176" = "Coneto De Comonfort"
    "B189" = "Nombre De Dios"
    "B193" = "Pánuco De Coronado"
    "B199" = "San Juan De Guadalupe"
    "B200" = "San Juan Del Río"
    "B201" = "San Pedro Del Gallo"
    "A211" = "Estado De México"
    "B211" = "Acambay De Ruíz Castañeda"
    "B214" = "Almoloya De Alquisiras"
    "B215" = "Almoloya De Juárez"
    "B221" = "Atizapán De Zaragoza"
    "B227" = "Chapa De Mota"
    "B230" = "Coacalco De Berriozábal"
    "B236" = "Ecatepec De Morelos"
    "B238" = "Ixtapan De La Sal"
    "B239" = "Ixtapan Del Oro"
    "B245" = "Naucalpan De Juárez"
    "B251" = "San Felipe Del Progreso"
    "B261" = "Tenango Del Valle"
    "B270" = "Tlalnepantla De Baz"
    "B275" = "Valle De Bravo"
    "B276" = "Valle De Chalco Solidaridad"
    "B277" = "Villa Del Carbón"
    "B287" = "San Miguel De Allende"
    "B288" = "Apaseo El Alto"
    "B289" = "Apaseo El Grande"
    "B297" = "Dolores Hidalgo Cuna De La Independencia Nacional"
    "B301" = "Jaral Del Progreso"
    "B309" = "Purísima Del Rincón"
    "B313" = "San Diego De La Unión"
    "B315" = "San Francisco Del Rincón"
    "B317" = "San Luis De La Paz"
    "B319" = "Silao De La Victoria"
    "B322" = "Valle De Santiago"
    "B328" = "Acapulco De Juárez"
    "B331" = "Ajuchitlán Del Progreso"
    "B332" = "Alcozauca De Guerrero"
    "B335" = "Atenango Del Río"
    "B337" = "Atoyac De Álvarez"
    "B338" = "Ayutla De Los Libres"
    "B341" = "Chilapa De Álvarez"
    "B342" = "Chilpancingo De Los Bravo"
    "B347" = "Coyuca De Benítez"
    "B348" = "Coyuca De Catalán"
    "B352" = "Cuetzala Del Progreso"
    "B353" = "Cutzamala De Pinzón"
    "B358" = "Huitzuco De Los Figueroa"
    "B359" = "Iguala De La Independencia"
    "B361" = "Ixcateopan De Cuauhtémoc"
    "B362" = "Zihuatanejo De Azueta"
    "B364" = "La Unión De Isidoro Montes De Oca"
    "B366" = "Mártir De Cuilapan"
    "B378" = "Taxco De Alarcón"
    "B380" = "Técpan De Galeana"
    "B382" = "Tepecoacuilco De Trujano"
    "B384" = "Tixtla De Guerrero"
    "B388" = "Tlalixtaquilla De Maldonado"
    "B389" = "Tlapa De Comonfort"
    "B399" = "Agua Blanca De Iturbide"
    "B403" = "Atotonilco El Grande"
    "B407" = "Cuautepec De Hinojosa"
    "B411" = "Huasca De Ocampo"
    "B413" = "Huejutla De Reyes"
    "B421" = "Mineral Del Monte"
    "B422" = "Mixquiahuala De Juárez"
    "B423" = "Molango De Escamilla"
    "B425" = "Omitlán De Juárez"
    "B426" = "Pachuca De Soto"
    "B428" = "Progreso De Obregón"
    "B434" = "Tenango De Doria"
    "B436" = "Tepeji Del Río De Ocampo"
    "B439" = "Tezontepec De Aldama"
    "B441" = "Tula De Allende"
    "B442" = "Tulancingo De Bravo"
    "B445" = "Zacualtipán De Ángeles"
    "B446" = "Zapotlán De Juárez"
    "B450" = "Acatlán De Juárez"
    "B451" = "Ahualulco De Mercado"
    "B457" = "Atotonilco El Alto"
    "B459" = "Autlán De Navarro"
    "B471" = "Cuautitlán De García Barragán"
    "B478" = "Encarnación De Díaz"
    "B484" = "Huejuquilla El Alto"
    "B485" = "Ixtlahuacán De Los Membrillos"
    "B486" = "Ixtlahuacán Del Río"
    "B490" = "Jilotlán De Los Dolores"
    "B496" = "La Manzanilla De La Paz"
    "B497" = "Lagos De Moreno"
    "B508" = "San Diego De Alejandría"
    "B510" = "San Juan De Los Lagos"
    "B512" = "San Martín De Bolaños"
    "B513" = "San Miguel El Alto"
    "B514" = "San Sebastián Del Oeste"
    "B515" = "Santa María De Los Ángeles"
    "B518" = "Talpa De Allende"
    "B519" = "Tamazula De Gordiano"
    "B524" = "Teocuitatlán De Corona"
    "B525" = "Tepatitlán De Morelos"
    "B527" = "Tizapán El Alto"
    "B528" = "Tlajomulco De Zúñiga"
    "B537" = "Valle De Guadalupe"
    "B542" = "Yahualica De González Gallo"
    "B543" = "Zacoalco De Torres"
    "B545" = "Zapotlán El Grande"
    "B565" = "Coalcomán De Vázquez Pallares"
    "B567" = "Cojumatlán De Régules"
    "B625" = "Tiquicheo De Nicolás Romero"
    "B645" = "Coatlán Del Río"
    "B655" = "Puente De Ixtla"
    "B660" = "Tetela Del Volcán"
    "B661" = "Tlaltizapán De Zapata"
    "B668" = "Zacualpan De Amilpas"
    "B672" = "Amatlán De Cañas"
    "B676" = "Ixtlán Del Río"
    "B683" = "Santa María Del Oro"
    "B696" = "San Nicolás De Los Garza"
    "B699" = "Acatlán De Pérez Figueroa"
    "B704" = "Chalcatongo De Hidalgo"
    "B706" = "Coicoyán De Las Flores"
    "B707" = "Cuilápam De Guerrero"
    "B708" = "Heroica Ciudad De Ejutla De Crespo"
    "B709" = "Heroica Ciudad De Huajuapan De León"
    "B710" = "Heroica Ciudad De Tlaxiaco"
    "B711" = "Ixtlán De Juárez"
    "B712" = "Heroica Ciudad De Juchitán De Zaragoza"
    "B720" = "Miahuatlán De Porfirio Díaz"
    "B721" = "Nejapa De Madero"
    "B722" = "Oaxaca De Juárez"
    "B723" = "Ocotlán De Morelos"
    "B724" = "Putla Villa De Guerrero"
    "B727" = "San Agustín De Las Juntas"
    "B757" = "San Mateo Del Mar"
    "B761" = "San Pablo Villa De Mitla"
    "B763" = "San Pedro El Alto"
    "B781" = "Santa María Del Tule"
    "B787" = "Santa María Jalapa Del Marqués"
    "B812" = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
    "B813" = "Tlacolula De Matamoros"
    "B814" = "Villa De Tututepec"
    "B815" = "Villa De Zaachila"
    "B816" = "Villa Sola De Vega"
    "B817" = "Villa Tejúpam De La Unión"
    "B818" = "Zimatlán De Álvarez"
    "B830" = "Chalchicomula De Sesma"
    "B841" = "Cuayuca De Andrade"
    "B849" = "Huehuetlán El Chico"
    "B853" = "Izúcar De Matamoros"
    "B860" = "Los Reyes De Juárez"
    "B861" = "Mazapiltepec De Juárez"
    "B864" = "Palmar De Bravo"
    "B872" = "San Nicolás De Los Ranchos"
    "B875" = "San Salvador El Seco"
    "B876" = "San Salvador El Verde"
    "B881" = "Tepanco De López"
    "B885" = "Tepexi De Rodríguez"
    "B886" = "Tetela De Ocampo"
    "B905" = "Amealco De Bonfil"
    "B907" = "Cadereyta De Montes"
    "B908" = "Jalpan De Serra"
    "B909" = "Landa De Matamoros"
    "B912" = "Pinal De Amoles"
    "B915" = "San Juan Del Río"
    "B921" = "Armadillo De Los Infante"
    "B924" = "Ciudad Del Maíz"
    "B933" = "San Ciro De Acosta"
    "B937" = "Santa María Del Río"
    "B944" = "Villa De Arriaga"
    "B945" = "Villa De Ramos"
    "B946" = "Villa De Reyes"
    "B1009" = "Nacozari De García"
    "B1021" = "San Miguel De Horcasitas"
    "B1022" = "San Pedro De La Cueva"
    "B1059" = "Soto La Marina"
    "B1070" = "Contla De Juan Cuamatzi"
    "B1073" = "Ixtacuixtla De Mariano Matamoros"
    "B1076" = "Papalotla De Xicohténcatl"
    "B1077" = "San Pablo Del Monte"
    "B1090" = "Amatlán De Los Reyes"
    "B1094" = "Boca Del Río"
    "B1096" = "Castillo De Teayo"
    "B1098" = "Cazones De Herrera"
    "B1105" = "Cosamaloapan De Carpio"
    "B1117" = "Hueyapan De Ocampo"
    "B1118" = "Ignacio De La Llave"
    "B1129" = "Las Vigas De Ramírez"
    "B1130" = "Lerdo De Tejada"
    "B1133" = "Martínez De La Torre"
    "B1147" = "Paso De Ovejas"
    "B1148" = "Paso Del Macho"
    "B1150" = "Poza Rica De Hidalgo"
    "B1158" = "Sayula De Alemán"
    "B1160" = "Soledad De Doblado"
    "B1183" = "Vega De Alatorre"
    "B1188" = "Zozocolco De Hidalgo"
    "B1201" = "Cañitas De Felipe Pescador"
    "B1203" = "Concepción Del Oro"
    "B1205" = "El Plateado De Joaquín Amaro"
    "B1215" = "Jiménez Del Teul"
    "B1221" = "Mezquital Del Oro"
    "B1225" = "Moyahua De Estrada"
    "B1226" = "Nochistlán De Mejía"
    "B1236" = "Teúl De González Ortega"
    "B1237" = "Tlaltenango De Sánchez Román"
    "B1240" = "Villa De Cos"
}
foreach ($cellRef in $textChanges.Keys) {
    $ws.Range($cellRef).Value = $textChanges[$cellRef]
}

# --- Fix minor floating point rounding (1-ULP) differences in percentage column ---
$valueChanges = @{
    "D19" = 0.0009477856281233844
    "D118" = 0.009133206961916248
    "D196" = 0.0009477856281233844
    "D291" = 0.0009477856281233844
    "D334" = 0.0009477856281233844
    "D358" = 0.0009477856281233844
    "D361" = 0.0009477856281233844
    "D403" = 0.0009477856281233844
    "D504" = 0.0009477856281233844
    "D546" = 0.0009477856281233844
    "D582" = 0.0009477856281233844
    "D589" = 0.0009477856281233844
    "D646" = 0.0009477856281233844
    "D679" = 0.0009477856281233844
    "D923" = 0.0009477856281233844
    "D1104" = 0.0009477856281233844
}
foreach ($cellRef in $valueChanges.Keys) {
    $ws.Range($cellRef).Value = $valueChanges[$cellRef]
}

# --- Remove trailing footer/metadata rows (1247:1251) ---
$ws.Rows("1247:1251").Delete()

Write-Host "Edit complete"